$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.485.35'
$ws.Range('E2').Value = '  +3.66%  '
$ws.Range('D3').Value = '2.482.82'
$ws.Range('E3').Value = '  +6.89%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '481.13'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +9.23%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '140.59'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +13.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.508'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +6.67%  '
$ws.Range('D9').Value = '2.489.69'
$ws.Range('E9').Value = '  +6.70%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0983'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.05%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.45'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.29%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.327'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +6.00%  '
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '2.921.16'
$ws.Range('E14').Value = '  +7.07%  '
$ws.Range('D15').Value = '55.507.29'
$ws.Range('E15').Value = '  +3.74%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.54'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +8.72%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000137'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +13.52%  '
$ws.Range('D18').Value = '2.490.81'
$ws.Range('E18').Value = '  +5.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.37'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +10.48%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '320.26'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +7.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.02'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +9.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.69'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '57.75'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.94%  '
$ws.Range('E25').Value = '  +10.07%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.408'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +10.85%  '
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('D28').Value = '2.605.20'
$ws.Range('E28').Value = '  +6.94%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.38'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +6.57%  '
$ws.Range('D30').Value = '0.0₃0789'
$ws.Range('E30').Value = '  +12.54%  '
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '148.87'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.15'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.48'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +9.56%  '
$ws.Range('E35').Value = '  +10.55%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.69'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +4.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.12'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +11.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.855'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.79%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '34.25'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.607'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +16.52%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0552'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +11.49%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.39'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +8.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.32'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +8.91%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.15'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('D46').Value = '1.971.97'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0902'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +8.68%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0223'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.67%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.46'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +9.77%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.53'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +14.03%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '248.24'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +29.97%  '
